$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$values = @{
    "H32" = 0
    "J32" = 0
    "L32" = 0
    "H42" = 1599
    "I42" = 2848
    "J42" = 350
    "K42" = 8544
    "L42" = 1050
    "M42" = -8314
    "N42" = -1510
    "H46" = 237833.17
    "J46" = 81750
    "L46" = 245250
    "N46" = -245488
    "H60" = 237833.17
    "J60" = 81750
    "L60" = 245250
    "N60" = -246218
    "H137" = 821.2
    "I137" = 821.2
    "K137" = 2463.6
    "M137" = 86.39999999999964
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
$clearRefs = @("N32")
foreach ($ref in $clearRefs) {
    $ws.Range($ref).ClearContents()
}

$ws = $wb.Worksheets.Item("ARM")
$values = @{
    "H2" = 1507.7391
    "I2" = 1768.7858
    "J2" = 1101.6666
    "K2" = 1768.7858
    "L2" = 1101.6666
    "M2" = -1655.7858
    "N2" = -1327.6666
    "H32" = 3655.3547
    "I32" = 3087.8704
    "K32" = 3087.8704
    "M32" = -2800.8704
    "H35" = 5369.6665
    "I35" = 5554.5
    "J35" = 5000
    "K35" = 5554.5
    "L35" = 5000
    "M35" = -5148.5
    "N35" = -5812
    "H45" = 7548.8237
    "I45" = 8238
    "J45" = 4332.6665
    "K45" = 8238
    "L45" = 4332.6665
    "M45" = -7861
    "N45" = -5086.6665
    "H88" = 1485.75
    "I88" = 1314.4
    "J88" = 1771.3334
    "K88" = 1314.4
    "L88" = 1771.3334
    "M88" = -908.4000000000001
    "N88" = -2583.3334
    "H91" = 1485.75
    "I91" = 1314.4
    "J91" = 1771.3334
    "K91" = 1314.4
    "L91" = 1771.3334
    "M91" = 89.59999999999991
    "N91" = -4579.3334
    "H116" = 1507.7391
    "I116" = 1768.7858
    "J116" = 1101.6666
    "K116" = 1768.7858
    "L116" = 1101.6666
    "M116" = 525.2141999999999
    "N116" = -5689.6666
    "H132" = 2309.0833
    "I132" = 2264.1875
    "J132" = 2668.25
    "K132" = 6792.5625
    "L132" = 8004.75
    "M132" = -4262.5625
    "N132" = -13064.75
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

$ws = $wb.Worksheets.Item("BSM")
$values = @{
    "H3" = 1507.7391
    "I3" = 1768.7858
    "J3" = 1101.6666
    "K3" = 1768.7858
    "L3" = 1101.6666
    "M3" = -1654.7858
    "N3" = -1329.6666
    "H22" = 998.6
    "I22" = 998
    "J22" = 999
    "K22" = 998
    "L22" = 999
    "M22" = -825
    "N22" = -1345
    "H86" = 30304868
    "I86" = 62501924
    "J86" = 1756.0588
    "K86" = 62501924
    "L86" = 1756.0588
    "M86" = -62500801
    "N86" = -4002.0588
    "H89" = 30304868
    "I89" = 62501924
    "J89" = 1756.0588
    "K89" = 312509620
    "L89" = 8780.294
    "M89" = -312504004
    "N89" = -20012.294
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

$ws = $wb.Worksheets.Item("CRP")
$values = @{
    "H62" = 345665
    "I62" = 1000000
    "J62" = 18497.5
    "K62" = 1000000
    "L62" = 18497.5
    "M62" = -999376
    "N62" = -19745.5
    "H65" = 345665
    "I65" = 1000000
    "J65" = 18497.5
    "K65" = 5000000
    "L65" = 92487.5
    "M65" = -4996880
    "N65" = -98727.5
    "H107" = 635.35486
    "J107" = 709.13336
    "L107" = 709.13336
    "N107" = -4549.13336
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

$ws = $wb.Worksheets.Item("CUL")
$values = @{
    "H11" = 113368.15
    "I11" = 475.9524
    "J11" = 350441.75
    "K11" = 1427.8572
    "L11" = 1051325.25
    "M11" = -1287.8572
    "N11" = -1051605.25
    "H81" = 4030
    "I81" = 1083.3334
    "J81" = 8450
    "K81" = 3250.0002
    "L81" = 25350
    "M81" = -2127.0002
    "N81" = -27596
    "H84" = 4030
    "I84" = 1083.3334
    "J84" = 8450
    "K84" = 9750.000599999999
    "L84" = 76050
    "M84" = -4134.000599999999
    "N84" = -87282
    "H121" = 11111965
    "J121" = 16667860
    "L121" = 50003580
    "N121" = -50006200
    "H126" = 10535.375
    "I126" = 2157.8
    "K126" = 6473.400000000001
    "M126" = -1533.400000000001
    "H127" = 57059
    "J127" = 64996
    "L127" = 194988
    "N127" = -204908
    "H129" = 2142.75
    "I129" = 2200
    "J129" = 2134.5715
    "K129" = 6600
    "L129" = 6403.7145
    "M129" = -1600
    "N129" = -16403.7145
    "H131" = 1268.1
    "I131" = 1268.1
    "J131" = 0
    "K131" = 3804.3
    "L131" = 0
    "M131" = 1235.7
    "H140" = 1970.125
    "I140" = 1663.2307
    "J140" = 3300
    "K140" = 4989.6921
    "L140" = 9900
    "M140" = 190.3078999999998
    "N140" = -20260
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
$clearRefs = @("N131")
foreach ($ref in $clearRefs) {
    $ws.Range($ref).ClearContents()
}

$ws = $wb.Worksheets.Item("GSM")
$values = @{
    "H2" = 160.08333
    "I2" = 24.666666
    "J2" = 566.3333
    "K2" = 24.666666
    "L2" = 566.3333
    "M2" = 88.33333400000001
    "N2" = -792.3333
    "H3" = 1004.5
    "I3" = 900
    "J3" = 1039.3334
    "K3" = 900
    "L3" = 1039.3334
    "M3" = -784
    "N3" = -1271.3334
    "H10" = 11216
    "I10" = 15500
    "J10" = 9502.4
    "K10" = 15500
    "L10" = 9502.4
    "M10" = -15331
    "N10" = -9840.4
    "H19" = 993.6667
    "I19" = 987.5
    "J19" = 1006
    "K19" = 987.5
    "L19" = 1006
    "M19" = -699.5
    "N19" = -1582
    "H32" = 0
    "J32" = 0
    "L32" = 0
    "H42" = 59290
    "J42" = 59290
    "L42" = 59290
    "N42" = -60260
    "H80" = 1548.5333
    "I80" = 1663.3
    "J80" = 1319
    "K80" = 1663.3
    "L80" = 1319
    "M80" = -665.3
    "N80" = -3315
    "H83" = 1548.5333
    "I83" = 1663.3
    "J83" = 1319
    "K83" = 8316.5
    "L83" = 6595
    "M83" = -3324.5
    "N83" = -16579
    "H97" = 75894590
    "I97" = 125001370
    "J97" = 2304.818
    "K97" = 125001370
    "L97" = 2304.818
    "M97" = -125000874
    "N97" = -3296.818
    "H115" = 59290
    "J115" = 59290
    "L115" = 59290
    "N115" = -61640
    "H126" = 26319262
    "I126" = 33335870
    "J126" = 6985
    "K126" = 100007610
    "L126" = 20955
    "M126" = -100005140
    "N126" = -25895
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
$clearRefs = @("N32")
foreach ($ref in $clearRefs) {
    $ws.Range($ref).ClearContents()
}

$ws = $wb.Worksheets.Item("WVR")
$values = @{
    "H27" = 0
    "J27" = 0
    "L27" = 0
    "H140" = 72107
    "J140" = 72107
    "L140" = 72107
    "N140" = -82467
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
$clearRefs = @("N27")
foreach ($ref in $clearRefs) {
    $ws.Range($ref).ClearContents()
}
